$d = $word.ActiveDocument

$d.Content.Find.Execute("363×6=2178", $true, $false, $false, $false, $false, $true, 1, $false, "180×5=900", 2) | Out-Null
$d.Content.Find.Execute("264×2=528", $true, $false, $false, $false, $false, $true, 1, $false, "138×2=276", 2) | Out-Null
$d.Content.Find.Execute("691×6=4146", $true, $false, $false, $false, $false, $true, 1, $false, "857×9=7713", 2) | Out-Null
$d.Content.Find.Execute("940×9=8460", $true, $false, $false, $false, $false, $true, 1, $false, "162×8=1296", 2) | Out-Null
$d.Content.Find.Execute("742×4=2968", $true, $false, $false, $false, $false, $true, 1, $false, "520×8=4160", 2) | Out-Null
$d.Content.Find.Execute("485×8=3880", $true, $false, $false, $false, $false, $true, 1, $false, "107×2=214", 2) | Out-Null
$d.Content.Find.Execute("845×8=6760", $true, $false, $false, $false, $false, $true, 1, $false, "979×2=1958", 2) | Out-Null
$d.Content.Find.Execute("660×5=3300", $true, $false, $false, $false, $false, $true, 1, $false, "188×2=376", 2) | Out-Null
$d.Content.Find.Execute("813×9=7317", $true, $false, $false, $false, $false, $true, 1, $false, "488×4=1952", 2) | Out-Null
$d.Content.Find.Execute("811×9=7299", $true, $false, $false, $false, $false, $true, 1, $false, "546×5=2730", 2) | Out-Null
$d.Content.Find.Execute("628×4=2512", $true, $false, $false, $false, $false, $true, 1, $false, "726×4=2904", 2) | Out-Null
$d.Content.Find.Execute("574×4=2296", $true, $false, $false, $false, $false, $true, 1, $false, "698×2=1396", 2) | Out-Null
$d.Content.Find.Execute("218×7=1526", $true, $false, $false, $false, $false, $true, 1, $false, "311×8=2488", 2) | Out-Null
$d.Content.Find.Execute("941×7=6587", $true, $false, $false, $false, $false, $true, 1, $false, "451×4=1804", 2) | Out-Null
$d.Content.Find.Execute("893×6=5358", $true, $false, $false, $false, $false, $true, 1, $false, "468×7=3276", 2) | Out-Null
$d.Content.Find.Execute("869×2=1738", $true, $false, $false, $false, $false, $true, 1, $false, "541×5=2705", 2) | Out-Null
$d.Content.Find.Execute("615×6=3690", $true, $false, $false, $false, $false, $true, 1, $false, "807×6=4842", 2) | Out-Null
$d.Content.Find.Execute("121×5=605", $true, $false, $false, $false, $false, $true, 1, $false, "252×6=1512", 2) | Out-Null
$d.Content.Find.Execute("300×2=600", $true, $false, $false, $false, $false, $true, 1, $false, "397×7=2779", 2) | Out-Null
$d.Content.Find.Execute("199×8=1592", $true, $false, $false, $false, $false, $true, 1, $false, "828×6=4968", 2) | Out-Null
$d.Content.Find.Execute("884×8=7072", $true, $false, $false, $false, $false, $true, 1, $false, "774×9=6966", 2) | Out-Null
$d.Content.Find.Execute("796×5=3980", $true, $false, $false, $false, $false, $true, 1, $false, "821×6=4926", 2) | Out-Null
$d.Content.Find.Execute("592×2=1184", $true, $false, $false, $false, $false, $true, 1, $false, "967×6=5802", 2) | Out-Null
$d.Content.Find.Execute("931×9=8379", $true, $false, $false, $false, $false, $true, 1, $false, "993×6=5958", 2) | Out-Null
$d.Content.Find.Execute("447×2=894", $true, $false, $false, $false, $false, $true, 1, $false, "792×5=3960", 2) | Out-Null
